# Add 4 new data rows (178-181) to the "dataset" worksheet, continuing the
# existing ersilia-os/ersilia entries, as described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 178

# Data for the new rows, in column order A..O
$rows = @(
    @("177", "https://github.com/ersilia-os/ersilia", "ersilia", "ersilia-os", "07/04/2020", "0", "0", "1", "1", "1", "0", "0", "0", "0", "1"),
    @("178", "https://github.com/ersilia-os/ersilia", "ersilia", "ersilia-os", "07/04/2020", "0", "0", "1", "1", "1", "0", "0", "0", "0", "1"),
    @("179", "https://github.com/ersilia-os/ersilia", "ersilia", "ersilia-os", "07/04/2020", "0", "0", "1", "1", "1", "0", "0", "0", "0", "1"),
    @("180", "https://github.com/ersilia-os/ersilia", "ersilia", "ersilia-os", "07/04/2020", "0", "1", "1", "1", "1", "0", "0", "0", "0", "1")
)

$endRow = $startRow + $rows.Count - 1

# Columns B..O in the existing sheet hold their data as literal text (not
# numbers/dates), matching the pattern used by every prior data row.
# Force the text number format first so values like "0"/"1"/"07/04/2020"
# are stored as text rather than being auto-converted to numbers/dates.
$ws.Range("B$($startRow):O$($endRow)").NumberFormat = "@"

# Column A keeps matching the bold/bordered/centered look used by every
# other data row (A4:A177) in the sheet.
$colA = $ws.Range("A$($startRow):A$($endRow)")
$colA.Font.Bold = $true
$colA.Borders.LineStyle = 1
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]

    $ws.Cells.Item($r, 1).Value2 = [int]$values[0]

    for ($c = 1; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $values[$c]
    }
}
